$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet is a weekly price log for "Locoto" at Vega Modelo de Temuco.
# A new week's record becomes the new row 12, and every existing record
# from the old row 12 through row 19 shifts down by one row (old row N ->
# new row N+1), with the last one ending up in the newly added row 20.
# Only the varying columns (D=Fecha, J=Volumen, K=Precio minimo,
# L=Precio maximo, M=Precio promedio ponderado, P=Precio $/Kg) differ
# between records; the rest (A, B, C, E, F, G, H, I, N, O, Q, R) are
# constant for every row in this sheet.

# 1) Build the new row 20 by copying every column from row 19 (the
#    constant columns already match what row 20 needs verbatim).
foreach ($col in @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R")) {
    $ws.Range("$col" + "20").Value = $ws.Range("$col" + "19").Value2
}
# Match the date-formatted style used by column D on the other data rows.
$ws.Range("D20").NumberFormat = $ws.Range("D19").NumberFormat

# 2) Cascade the varying data down one row at a time, starting from the
#    bottom so no value is overwritten before it has been copied onward.
for ($r = 19; $r -ge 13; $r--) {
    $src = $r - 1
    $ws.Range("D$r").Value = $ws.Range("D$src").Value2
    $ws.Range("J$r").Value = $ws.Range("J$src").Value2
    $ws.Range("K$r").Value = $ws.Range("K$src").Value2
    $ws.Range("L$r").Value = $ws.Range("L$src").Value2
    $ws.Range("M$r").Value = $ws.Range("M$src").Value2
    $ws.Range("P$r").Value = $ws.Range("P$src").Value2
}

# 3) Write the new week's record into row 12.
$ws.Range("D12").Value = 44719
$ws.Range("J12").Value = 80
$ws.Range("K12").Value = 3600
$ws.Range("L12").Value = 3600
$ws.Range("M12").Value = 3600
$ws.Range("P12").Value = 3600
